# Insert a new "date" column at the front of the data rows (rows 2-51),
# shifting the existing columns A:L -> B:M for each of those rows, and
# filling the new column A with the literal text "2023/02/16".
# Row 1 (the header row) is left untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 51
$lastCol = 12   # column L in the original layout

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    # Walk from the right-most column down to column A, copying each
    # cell's value one column to the right so nothing is clobbered.
    for ($c = $lastCol; $c -ge 1; $c--) {
        $srcVal = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($r, $c + 1).Value = $srcVal
    }
    # New column A gets the literal date string (kept as text).
    $ws.Cells.Item($r, 1).Value = "2023/02/16"
}

# The source diff shows row 7's (new) category cell (D7) truncated to
# its first 35 characters - replicate that exact edit.
$fullCategory = $ws.Cells.Item(7, 4).Value()
$ws.Cells.Item(7, 4).Value = $fullCategory.Substring(0, 35)

# Sheet's used range grew by one column (A1:L51 -> A1:M51).
$ws.Range("A1:M51").Select()
